# Refresh the crypto price/volume table (and the two pairs of rows whose
# ranking order swapped) to match the latest scrape.
#
# Columns D (Price) and E (Volume(1h)) are stored as *text* in the sheet
# (e.g. "1.00", "0.999", "  +0.95%  ") rather than numbers, so several of
# the new Price values look numeric ("1.00", "2.00", "0.999", ...). Those
# are written with a leading apostrophe -- Excel's standard "force text"
# prefix -- so they keep their exact text form instead of being
# auto-coerced into numbers (which would silently drop trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.367.05'
$ws.Range("E2").Value = '  +0.95%  '
# Row 3
$ws.Range("D3").Value = '2.570.11'
$ws.Range("E3").Value = '  +1.28%  '
# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.09%  '
# Row 5
$ws.Range("D5").Value = '''585.21'
$ws.Range("E5").Value = '  +3.30%  '
# Row 6
$ws.Range("D6").Value = '''148.35'
$ws.Range("E6").Value = '  +1.17%  '
# Row 7
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.07%  '
# Row 8
$ws.Range("D8").Value = '''0.604'
$ws.Range("E8").Value = '  +4.50%  '
# Row 9
$ws.Range("E9").Value = '  +4.46%  '
# Row 10
$ws.Range("D10").Value = '''5.68'
$ws.Range("E10").Value = '  +1.58%  '
# Row 11
$ws.Range("E11").Value = '  +0.48%  '
# Row 12
$ws.Range("E12").Value = '  +1.85%  '
# Row 13
$ws.Range("D13").Value = '''27.54'
$ws.Range("E13").Value = '  +2.38%  '
# Row 14
$ws.Range("D14").Value = '3.030.11'
$ws.Range("E14").Value = '  +1.21%  '
# Row 15
$ws.Range("D15").Value = '63.267.06'
$ws.Range("E15").Value = '  +0.77%  '
# Row 16
$ws.Range("E16").Value = '  +5.54%  '
# Row 17
$ws.Range("D17").Value = '2.584.40'
$ws.Range("E17").Value = '  +1.72%  '
# Row 18
$ws.Range("D18").Value = '''11.39'
$ws.Range("E18").Value = '  -0.44%  '
# Row 19
$ws.Range("D19").Value = '''343.64'
$ws.Range("E19").Value = '  +3.04%  '
# Row 20
$ws.Range("E20").Value = '  +3.79%  '
# Row 21
$ws.Range("D21").Value = '''6.89'
$ws.Range("E21").Value = '  +1.95%  '
# Row 22
$ws.Range("E22").Value = '  +0.18%  '
# Row 23
$ws.Range("D23").Value = '''66.76'
$ws.Range("E23").Value = '  +3.25%  '
# Row 24
$ws.Range("B24").Value = 'Fetch.AI'
$ws.Range("C24").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D24").Value = '''1.65'
$ws.Range("E24").Value = '  +4.18%  '
# Row 25
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '2.686.53'
$ws.Range("E25").Value = '  +0.70%  '
# Row 26
$ws.Range("E26").Value = '  +0.81%  '
# Row 27
$ws.Range("D27").Value = '''8.28'
$ws.Range("E27").Value = '  +14.57%  '
# Row 28
$ws.Range("D28").Value = '''8.56'
$ws.Range("E28").Value = '  +3.09%  '
# Row 29
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.08%  '
# Row 30
$ws.Range("E30").Value = '  +0.57%  '
# Row 31
$ws.Range("D31").Value = '''2.00'
$ws.Range("E31").Value = '  +8.51%  '
# Row 32
$ws.Range("D32").Value = '0.0₃0830'
$ws.Range("E32").Value = '  +2.87%  '
# Row 33
$ws.Range("D33").Value = '''465.90'
$ws.Range("E33").Value = '  +15.27%  '
# Row 34
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''1.63'
$ws.Range("E34").Value = '  +4.13%  '
# Row 35
$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").Value = '''176.73'
$ws.Range("E35").Value = '  -0.01%  '
# Row 36
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '''19.31'
$ws.Range("E36").Value = '  +2.21%  '
# Row 37
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = '''0.403'
$ws.Range("E37").Value = '  +2.20%  '
# Row 38
$ws.Range("D38").Value = '''4.53'
$ws.Range("E38").Value = '  +5.14%  '
# Row 41
$ws.Range("E41").Value = '  -0.04%  '
# Row 42
$ws.Range("D42").Value = '''152.03'
$ws.Range("E42").Value = '  +0.37%  '
# Row 43
$ws.Range("E43").Value = '  +2.61%  '
# Row 44
$ws.Range("D44").Value = '''21.23'
$ws.Range("E44").Value = '  +3.43%  '
# Row 45
$ws.Range("E45").Value = '  +7.42%  '
# Row 46
$ws.Range("D46").Value = '''0.617'
$ws.Range("E46").Value = '  +2.79%  '
# Row 47
$ws.Range("D47").Value = '''0.0983'
$ws.Range("E47").Value = '  +2.91%  '
# Row 48
$ws.Range("E48").Value = '  +2.39%  '
# Row 49
$ws.Range("D49").Value = '''18.54'
$ws.Range("E49").Value = '  +1.84%  '
# Row 50
$ws.Range("D50").Value = '''1.76'
$ws.Range("E50").Value = '  +0.34%  '
# Row 51
$ws.Range("D51").Value = '''11.38'
$ws.Range("E51").Value = '  -0.15%  '
